$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.530.47'
$ws.Range("E2").Value = '  -2.66%  '
$ws.Range("D3").Value = '1.813.67'
$ws.Range("E3").Value = '  -2.35%  '
$ws.Range("D4").Value = "'" + '1.008'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.67%  '
$ws.Range("B5").Value = 'USDC'
$ws.Range("C5").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D5").Value = "'" + '1.007'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.52%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = "'" + '308.42'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.87%  '
$ws.Range("D7").Value = "'" + '0.4566'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.91%  '
$ws.Range("D8").Value = "'" + '0.3667'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.19%  '
$ws.Range("E9").Value = '  -2.33%  '
$ws.Range("D10").Value = "'" + '0.8798'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.39%  '
$ws.Range("D11").Value = "'" + '0.07748'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.62%  '
$ws.Range("D12").Value = "'" + '19.35'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.59%  '
$ws.Range("D13").Value = '1.833.86'
$ws.Range("E13").Value = '  -0.85%  '
$ws.Range("D14").Value = "'" + '5.295'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.20%  '
$ws.Range("D15").Value = "'" + '6.374'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'" + '86.67'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.39%  '
$ws.Range("D17").Value = "'" + '1.009'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.66%  '
$ws.Range("D18").Value = "'" + '0.000008603'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.74%  '
$ws.Range("E19").Value = '  +0.59%  '
$ws.Range("D20").Value = '26.593.65'
$ws.Range("E20").Value = '  -2.54%  '
$ws.Range("E21").Value = '  -3.35%  '
$ws.Range("D22").Value = "'" + '5.013'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.46%  '
$ws.Range("D23").Value = "'" + '10.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.47%  '
$ws.Range("D24").Value = "'" + '1.976'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.76%  '
$ws.Range("D25").Value = "'" + '151.36'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.28%  '
$ws.Range("E26").Value = '  -2.67%  '
$ws.Range("D27").Value = "'" + '2.073'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.15%  '
$ws.Range("D28").Value = "'" + '112.97'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.60%  '
$ws.Range("D29").Value = "'" + '4.852'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.93%  '
$ws.Range("D30").Value = "'" + '0.08701'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.63%  '
$ws.Range("E31").Value = '  -3.38%  '
$ws.Range("D32").Value = "'" + '4.518'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.25%  '
$ws.Range("E33").Value = '  -4.79%  '
$ws.Range("D34").Value = "'" + '2.693'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.14%  '
$ws.Range("E35").Value = '  -4.31%  '
$ws.Range("E36").Value = '  +0.93%  '
$ws.Range("D37").Value = "'" + '1.085'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.06%  '
$ws.Range("E38").Value = '  +0.68%  '
$ws.Range("D39").Value = "'" + '0.05116'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.20%  '
$ws.Range("D40").Value = "'" + '2.894'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.77%  '
$ws.Range("D41").Value = "'" + '7.007'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.79%  '
$ws.Range("D42").Value = "'" + '0.5001'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.28%  '
$ws.Range("D43").Value = "'" + '0.1554'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.63%  '
$ws.Range("D44").Value = "'" + '8.164'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.36%  '
$ws.Range("D45").Value = "'" + '1.007'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.54%  '
$ws.Range("D46").Value = "'" + '0.4608'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.05%  '
$ws.Range("D47").Value = "'" + '10.02'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.32%  '
$ws.Range("D48").Value = "'" + '101.38'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.35%  '
$ws.Range("E49").Value = '  -3.42%  '
$ws.Range("D50").Value = "'" + '0.06001'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.18%  '
$ws.Range("D51").Value = "'" + '64.53'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.51%  '
